$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.507950666666667
$ws.Range("H2").Value = 4.523852
$ws.Range("I2").Value = 0.2206625213859172
$ws.Range("J2").Value = 0.2310732482744153
$ws.Range("M2").Value = 1.037532
$ws.Range("N2").Value = 3.112596
$ws.Range("O2").Value = 0.04166450179684251
$ws.Range("P2").Value = 0.0439159257402554
$ws.Range("Q2").Value = 1.564547071088
$ws.Range("R2").Value = 14.080923639792
$ws.Range("S2").Value = 0.009193794018779344
$ws.Range("T2").Value = 0.01014779561177882
$ws.Range("G3").Value = 1.507950666666667
$ws.Range("H3").Value = 4.523852
$ws.Range("I3").Value = 0.2206625213859172
$ws.Range("J3").Value = 0.2310732482744153
$ws.Range("O3").Value = 0.2561129158441639
$ws.Range("P3").Value = 0.2699524849277078
$ws.Range("Q3").Value = 9.617316782176443
$ws.Range("R3").Value = 86.55585103958798
$ws.Range("S3").Value = 0.05651452176967243
$ws.Range("T3").Value = 0.06237879757199558
$ws.Range("G4").Value = 1.507950666666667
$ws.Range("H4").Value = 4.523852
$ws.Range("I4").Value = 0.2206625213859172
$ws.Range("J4").Value = 0.2310732482744153
$ws.Range("M4").Value = 6.239319333333333
$ws.Range("N4").Value = 18.717958
$ws.Range("O4").Value = 0.2505543265891952
$ws.Range("P4").Value = 0.2640935262839185
$ws.Range("Q4").Value = 9.408585748246223
$ws.Range("R4").Value = 84.67727173421599
$ws.Range("S4").Value = 0.05528794944932236
$ws.Range("T4").Value = 0.06102494896666972
$ws.Range("G5").Value = 1.507950666666667
$ws.Range("H5").Value = 4.523852
$ws.Range("I5").Value = 0.2206625213859172
$ws.Range("J5").Value = 0.2310732482744153
$ws.Range("M5").Value = 3.8299385
$ws.Range("N5").Value = 7.659877
$ws.Range("O5").Value = 0.1538000558200097
$ws.Range("P5").Value = 0.1080739644693659
$ws.Range("Q5").Value = 5.775358314367333
$ws.Range("R5").Value = 34.652149886204
$ws.Range("S5").Value = 0.03393790810653814
$ws.Range("T5").Value = 0.02497300202383012
$ws.Range("G6").Value = 1.507950666666667
$ws.Range("H6").Value = 4.523852
$ws.Range("I6").Value = 0.2206625213859172
$ws.Range("J6").Value = 0.2310732482744153
$ws.Range("M6").Value = 7.417532333333334
$ws.Range("N6").Value = 22.252597
$ws.Range("O6").Value = 0.2978681999497886
$ws.Range("P6").Value = 0.3139640985787523
$ws.Range("Q6").Value = 11.18527282707156
$ws.Range("R6").Value = 100.667455443644
$ws.Range("S6").Value = 0.06572834804160489
$ws.Range("T6").Value = 0.07254870410014104
$ws.Range("I7").Value = 0.323871134529766
$ws.Range("J7").Value = 0.3391511825754447
$ws.Range("M7").Value = 1.037532
$ws.Range("N7").Value = 3.112596
$ws.Range("O7").Value = 0.04166450179684251
$ws.Range("P7").Value = 0.0439159257402554
$ws.Range("Q7").Value = 2.29631942822
$ws.Range("R7").Value = 20.66687485398
$ws.Range("S7").Value = 0.01349392946656086
$ws.Range("T7").Value = 0.01489413814870303
$ws.Range("I8").Value = 0.323871134529766
$ws.Range("J8").Value = 0.3391511825754447
$ws.Range("O8").Value = 0.2561129158441639
$ws.Range("P8").Value = 0.2699524849277078
$ws.Range("S8").Value = 0.08294758062217583
$ws.Range("T8").Value = 0.091554704502412
$ws.Range("I9").Value = 0.323871134529766
$ws.Range("J9").Value = 0.3391511825754447
$ws.Range("M9").Value = 6.239319333333333
$ws.Range("N9").Value = 18.717958
$ws.Range("O9").Value = 0.2505543265891952
$ws.Range("P9").Value = 0.2640935262839185
$ws.Range("Q9").Value = 13.80918391336556
$ws.Range("R9").Value = 124.28265522029
$ws.Range("S9").Value = 0.08114731401378415
$ws.Range("T9").Value = 0.08956763174971023
$ws.Range("I10").Value = 0.323871134529766
$ws.Range("J10").Value = 0.3391511825754447
$ws.Range("M10").Value = 3.8299385
$ws.Range("N10").Value = 7.659877
$ws.Range("O10").Value = 0.1538000558200097
$ws.Range("P10").Value = 0.1080739644693659
$ws.Range("Q10").Value = 8.476617768355833
$ws.Range("R10").Value = 50.859706610135
$ws.Range("S10").Value = 0.04981139856916787
$ws.Range("T10").Value = 0.03665341285540202
$ws.Range("I11").Value = 0.323871134529766
$ws.Range("J11").Value = 0.3391511825754447
$ws.Range("M11").Value = 7.417532333333334
$ws.Range("N11").Value = 22.252597
$ws.Range("O11").Value = 0.2978681999497886
$ws.Range("P11").Value = 0.3139640985787523
$ws.Range("Q11").Value = 16.41686579930389
$ws.Range("R11").Value = 147.751792193735
$ws.Range("S11").Value = 0.09647091185807723
$ws.Range("T11").Value = 0.1064812953192174
$ws.Range("G12").Value = 1.299855666666667
$ws.Range("H12").Value = 3.899567
$ws.Range("I12").Value = 0.190211414195981
$ws.Range("J12").Value = 0.1991854759071952
$ws.Range("M12").Value = 1.037532
$ws.Range("N12").Value = 3.112596
$ws.Range("O12").Value = 0.04166450179684251
$ws.Range("P12").Value = 0.0439159257402554
$ws.Range("Q12").Value = 1.348641849548
$ws.Range("R12").Value = 12.137776645932
$ws.Range("S12").Value = 0.007925063808548403
$ws.Range("T12").Value = 0.008747414568477817
$ws.Range("G13").Value = 1.299855666666667
$ws.Range("H13").Value = 3.899567
$ws.Range("I13").Value = 0.190211414195981
$ws.Range("J13").Value = 0.1991854759071952
$ws.Range("O13").Value = 0.2561129158441639
$ws.Range("P13").Value = 0.2699524849277078
$ws.Range("Q13").Value = 8.290141046241445
$ws.Range("R13").Value = 74.61126941617299
$ws.Range("S13").Value = 0.04871559991657468
$ws.Range("T13").Value = 0.05377061418265543
$ws.Range("G14").Value = 1.299855666666667
$ws.Range("H14").Value = 3.899567
$ws.Range("I14").Value = 0.190211414195981
$ws.Range("J14").Value = 0.1991854759071952
$ws.Range("M14").Value = 6.239319333333333
$ws.Range("N14").Value = 18.717958
$ws.Range("O14").Value = 0.2505543265891952
$ws.Range("P14").Value = 0.2640935262839185
$ws.Range("Q14").Value = 8.110214591576224
$ws.Range("R14").Value = 72.991931324186
$ws.Range("S14").Value = 0.0476582927934525
$ws.Range("T14").Value = 0.05260359471687168
$ws.Range("G15").Value = 1.299855666666667
$ws.Range("H15").Value = 3.899567
$ws.Range("I15").Value = 0.190211414195981
$ws.Range("J15").Value = 0.1991854759071952
$ws.Range("M15").Value = 3.8299385
$ws.Range("N15").Value = 7.659877
$ws.Range("O15").Value = 0.1538000558200097
$ws.Range("P15").Value = 0.1080739644693659
$ws.Range("Q15").Value = 4.978367262209834
$ws.Range("R15").Value = 29.870203573259
$ws.Range("S15").Value = 0.02925452612094486
$ws.Range("T15").Value = 0.02152676404600795
$ws.Range("G16").Value = 1.299855666666667
$ws.Range("H16").Value = 3.899567
$ws.Range("I16").Value = 0.190211414195981
$ws.Range("J16").Value = 0.1991854759071952
$ws.Range("M16").Value = 7.417532333333334
$ws.Range("N16").Value = 22.252597
$ws.Range("O16").Value = 0.2978681999497886
$ws.Range("P16").Value = 0.3139640985787523
$ws.Range("Q16").Value = 9.641721436166558
$ws.Range("R16").Value = 86.77549292549901
$ws.Range("S16").Value = 0.05665793155646053
$ws.Range("T16").Value = 0.06253708839318235
$ws.Range("G17").Value = 0.923658
$ws.Range("H17").Value = 1.847316
$ws.Range("I17").Value = 0.1351613866976242
$ws.Range("J17").Value = 0.09435881384035105
$ws.Range("M17").Value = 1.037532
$ws.Range("N17").Value = 3.112596
$ws.Range("O17").Value = 0.04166450179684251
$ws.Range("P17").Value = 0.0439159257402554
$ws.Range("Q17").Value = 0.9583247320559999
$ws.Range("R17").Value = 5.749948392336
$ws.Range("S17").Value = 0.005631431838926887
$ws.Range("T17").Value = 0.00414385466155144
$ws.Range("G18").Value = 0.923658
$ws.Range("H18").Value = 1.847316
$ws.Range("I18").Value = 0.1351613866976242
$ws.Range("J18").Value = 0.09435881384035105
$ws.Range("O18").Value = 0.2561129158441639
$ws.Range("P18").Value = 0.2699524849277078
$ws.Range("Q18").Value = 5.890850265034
$ws.Range("R18").Value = 35.345101590204
$ws.Range("S18").Value = 0.03461657685666911
$ws.Range("T18").Value = 0.02547239627103375
$ws.Range("G19").Value = 0.923658
$ws.Range("H19").Value = 1.847316
$ws.Range("I19").Value = 0.1351613866976242
$ws.Range("J19").Value = 0.09435881384035105
$ws.Range("M19").Value = 6.239319333333333
$ws.Range("N19").Value = 18.717958
$ws.Range("O19").Value = 0.2505543265891952
$ws.Range("P19").Value = 0.2640935262839185
$ws.Range("Q19").Value = 5.762997216788
$ws.Range("R19").Value = 34.577983300728
$ws.Range("S19").Value = 0.03386527022488503
$ws.Range("T19").Value = 0.02491955188306612
$ws.Range("G20").Value = 0.923658
$ws.Range("H20").Value = 1.847316
$ws.Range("I20").Value = 0.1351613866976242
$ws.Range("J20").Value = 0.09435881384035105
$ws.Range("M20").Value = 3.8299385
$ws.Range("N20").Value = 7.659877
$ws.Range("O20").Value = 0.1538000558200097
$ws.Range("P20").Value = 0.1080739644693659
$ws.Range("Q20").Value = 3.537553335033
$ws.Range("R20").Value = 14.150213340132
$ws.Range("S20").Value = 0.02078782881880452
$ws.Range("T20").Value = 0.01019773109435361
$ws.Range("G21").Value = 0.923658
$ws.Range("H21").Value = 1.847316
$ws.Range("I21").Value = 0.1351613866976242
$ws.Range("J21").Value = 0.09435881384035105
$ws.Range("M21").Value = 7.417532333333334
$ws.Range("N21").Value = 22.252597
$ws.Range("O21").Value = 0.2978681999497886
$ws.Range("P21").Value = 0.3139640985787523
$ws.Range("Q21").Value = 6.851263079942
$ws.Range("R21").Value = 41.10757847965201
$ws.Range("S21").Value = 0.04026027895833862
$ws.Range("T21").Value = 0.02962527993034612
$ws.Range("G22").Value = 0.8890256666666666
$ws.Range("H22").Value = 2.667077
$ws.Range("I22").Value = 0.1300935431907118
$ws.Range("J22").Value = 0.1362312794025938
$ws.Range("M22").Value = 1.037532
$ws.Range("N22").Value = 3.112596
$ws.Range("O22").Value = 0.04166450179684251
$ws.Range("P22").Value = 0.0439159257402554
$ws.Range("Q22").Value = 0.9223925779879999
$ws.Range("R22").Value = 8.301533201891999
$ws.Range("S22").Value = 0.005420282664027018
$ws.Range("T22").Value = 0.005982722749744295
$ws.Range("G23").Value = 0.8890256666666666
$ws.Range("H23").Value = 2.667077
$ws.Range("I23").Value = 0.1300935431907118
$ws.Range("J23").Value = 0.1362312794025938
$ws.Range("O23").Value = 0.2561129158441639
$ws.Range("P23").Value = 0.2699524849277078
$ws.Range("Q23").Value = 5.669974258984777
$ws.Range("R23").Value = 51.02976833086299
$ws.Range("S23").Value = 0.03331863667907186
$ws.Range("T23").Value = 0.03677597239961105
$ws.Range("G24").Value = 0.8890256666666666
$ws.Range("H24").Value = 2.667077
$ws.Range("I24").Value = 0.1300935431907118
$ws.Range("J24").Value = 0.1362312794025938
$ws.Range("M24").Value = 6.239319333333333
$ws.Range("N24").Value = 18.717958
$ws.Range("O24").Value = 0.2505543265891952
$ws.Range("P24").Value = 0.2640935262839185
$ws.Range("Q24").Value = 5.546915029862888
$ws.Range("R24").Value = 49.92223526876599
$ws.Range("S24").Value = 0.03259550010775116
$ws.Range("T24").Value = 0.03597779896760075
$ws.Range("G25").Value = 0.8890256666666666
$ws.Range("H25").Value = 2.667077
$ws.Range("I25").Value = 0.1300935431907118
$ws.Range("J25").Value = 0.1362312794025938
$ws.Range("M25").Value = 3.8299385
$ws.Range("N25").Value = 7.659877
$ws.Range("O25").Value = 0.1538000558200097
$ws.Range("P25").Value = 0.1080739644693659
$ws.Range("Q25").Value = 3.404913628254833
$ws.Range("R25").Value = 20.429481769529
$ws.Range("S25").Value = 0.02000839420455431
$ws.Range("T25").Value = 0.01472305444977218
$ws.Range("G26").Value = 0.8890256666666666
$ws.Range("H26").Value = 2.667077
$ws.Range("I26").Value = 0.1300935431907118
$ws.Range("J26").Value = 0.1362312794025938
$ws.Range("M26").Value = 7.417532333333334
$ws.Range("N26").Value = 22.252597
$ws.Range("O26").Value = 0.2978681999497886
$ws.Range("P26").Value = 0.3139640985787523
$ws.Range("Q26").Value = 6.594376627663222
$ws.Range("R26").Value = 59.349389648969
$ws.Range("S26").Value = 0.03875072953530739
$ws.Range("T26").Value = 0.04277173083586552
